# Refresh the cryptocurrency price/volume snapshot (D = Price, E = Volume(1h))
# for rows 2..51 of the active sheet, matching the latest scrape.
#
# D-column values are plain decimal-looking strings (e.g. "2.05", "1.00")
# that Excel's normal text-entry coercion would silently reinterpret as
# numbers, losing the original "number formatted as text" representation
# used throughout the sheet. To keep them as literal text (matching every
# other cell in the column) we round-trip each value through a text
# formula ( ="value" ) and then Copy/PasteSpecial values-only on top of
# itself — this bakes in a plain string without touching NumberFormat/
# style of the cell.
#
# E-column values already contain padding spaces and a trailing "%" sign,
# so Excel's COM Value setter naturally stores them as text and a direct
# assignment is sufficient.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D="63.787.71"; E="  +0.27%  " }
    @{ Row=3; D="2.631.09"; E="  +0.52%  " }
    @{ Row=4; D=$null; E="  -0.04%  " }
    @{ Row=5; D="578.26"; E="  +0.76%  " }
    @{ Row=6; D="156.70"; E="  +1.13%  " }
    @{ Row=7; D=$null; E="  +0.19%  " }
    @{ Row=8; D=$null; E="  -0.02%  " }
    @{ Row=9; D=$null; E="  -1.65%  " }
    @{ Row=10; D=$null; E="  +0.35%  " }
    @{ Row=11; D="0.385"; E="  +0.07%  " }
    @{ Row=12; D=$null; E="  +0.79%  " }
    @{ Row=13; D="28.67"; E="  +1.82%  " }
    @{ Row=14; D="3.105.64"; E="  +0.44%  " }
    @{ Row=15; D=$null; E="  +0.84%  " }
    @{ Row=16; D="63.700.70"; E="  +0.25%  " }
    @{ Row=17; D="2.657.49"; E="  +0.70%  " }
    @{ Row=18; D="12.16"; E="  +0.96%  " }
    @{ Row=19; D=$null; E="  +2.75%  " }
    @{ Row=20; D=$null; E="  -2.25%  " }
    @{ Row=21; D="343.43"; E="  +0.01%  " }
    @{ Row=22; D=$null; E="  -0.03%  " }
    @{ Row=23; D="68.30"; E="  +1.80%  " }
    @{ Row=24; D=$null; E="  +9.18%  " }
    @{ Row=25; D=$null; E="  +3.24%  " }
    @{ Row=26; D=$null; E="  +4.29%  " }
    @{ Row=27; D="9.23"; E="  +0.33%  " }
    @{ Row=28; D="581.58"; E="  -0.24%  " }
    @{ Row=29; D="8.26"; E="  +4.98%  " }
    @{ Row=30; D="1.00"; E="  +0.46%  " }
    @{ Row=31; D=$null; E="  -0.23%  " }
    @{ Row=32; D="2.05"; E="  -0.62%  " }
    @{ Row=33; D="1.73"; E="  +2.23%  " }
    @{ Row=34; D="6.64"; E="  +2.51%  " }
    @{ Row=35; D="5.46"; E="  +2.99%  " }
    @{ Row=36; D=$null; E="  -1.16%  " }
    @{ Row=37; D=$null; E="  -0.32%  " }
    @{ Row=38; D=$null; E="  +0.02%  " }
    @{ Row=39; D=$null; E="  +2.57%  " }
    @{ Row=40; D="153.78"; E="  +0.12%  " }
    @{ Row=41; D=$null; E="  +8.40%  " }
    @{ Row=42; D=$null; E="  -0.01%  " }
    @{ Row=43; D="163.31"; E="  +4.33%  " }
    @{ Row=44; D="24.06"; E="  +5.80%  " }
    @{ Row=45; D=$null; E="  -0.24%  " }
    @{ Row=46; D="0.0586"; E="  -1.31%  " }
    @{ Row=47; D=$null; E="  +0.68%  " }
    @{ Row=48; D=$null; E="  -1.10%  " }
    @{ Row=49; D=$null; E="  -0.96%  " }
    @{ Row=50; D=$null; E="  +1.65%  " }
    @{ Row=51; D=$null; E="  +2.12%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cell = $ws.Range("D$row")
        $cell.Formula = '="' + $u.D + '"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)
    }

    if ($null -ne $u.E) {
        $ws.Range("E$row").Value = $u.E
    }
}
